$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1167
$ws1.Range("F3").Value = 1080
$ws1.Range("F4").Value = 1878
$ws1.Range("F7").Value = 66
$ws1.Range("F8").Value = 22
$ws1.Range("F9").Value = 132
$ws1.Range("F10").Value = 323
$ws1.Range("F11").Value = 103
$ws1.Range("F13").Value = 770
$ws1.Range("F19").Value = 194
$ws1.Range("F21").Value = 61
$ws1.Range("F24").Value = 43
$ws1.Range("F25").Value = 895
$ws1.Range("F26").Value = 338
$ws1.Range("F27").Value = 181

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 19
$ws2.Range("F6").Value = 29
$ws2.Range("F8").Value = 83
$ws2.Range("F11").Value = 126

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 322

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 322
$ws4.Range("F3").Value = 1167
$ws4.Range("F4").Value = 1080
$ws4.Range("F5").Value = 1878
$ws4.Range("F8").Value = 66
$ws4.Range("F10").Value = 22
$ws4.Range("F11").Value = 132
$ws4.Range("F12").Value = 323
$ws4.Range("F13").Value = 103
$ws4.Range("F15").Value = 770
$ws4.Range("F22").Value = 19
$ws4.Range("F24").Value = 29
$ws4.Range("F27").Value = 194
$ws4.Range("F29").Value = 61
$ws4.Range("F32").Value = 43
$ws4.Range("F33").Value = 895
$ws4.Range("F34").Value = 338
$ws4.Range("F35").Value = 83
$ws4.Range("F37").Value = 181
$ws4.Range("F41").Value = 126
$ws4.Range("F42").Value = 126

$wb.Save()
